$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 6 (-25C) data
$ws.Range("B6").Value = 9925
$ws.Range("C6").Value = 4.1036000000000001
$ws.Range("D6").Value = 3.5800999999999998
$ws.Range("E6").Value = -0.00107
$ws.Range("F6").Value = -3.9947499999999998
$ws.Range("H6").Value = 69.900000000000006

# Fill in row 7 (-30C) data
$ws.Range("B7").Value = 10145
$ws.Range("C7").Value = 4.0776000000000003
$ws.Range("D7").Value = 3.4533499999999999
$ws.Range("E7").Value = -0.00106
$ws.Range("F7").Value = -3.9948000000000001
$ws.Range("H7").Value = 65.7

# Recalculate so G6/G7 formulas resolve (they are shared formulas already present)
$excel.Calculate()

# Update the selection to reflect where the user left off (G7)
$ws.Range("G7").Select()
